# Release 2024-03: Excel files ready for translation
#
# slugstatus.farming.co.uk.xlsx translation workbook:
#   "main" sheet          - row 2: product name (KEY slugstatus.farming.co.uk.1_0.name)
#   "SlugWatch2023" sheet - row 2: model description, row 3: model name, row 4: model purpose
#
# Columns C..M are the newly-delivered translations (dk, de, nl, sl, no, se, fi, lt, fr, it, gr).
# "SlugWatch" (the product/model name) is a brand name and stays the same in every language;
# the description/purpose text is the actual translated copy for each locale.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "main"
$ws2 = $wb.Worksheets.Item(2)   # "SlugWatch2023"

# --- main!row2 : product name ("SlugWatch") is unchanged across languages ---
$ws1.Range("C2").Value = "SlugWatch"
$ws1.Range("D2").Value = "SlugWatch"
$ws1.Range("F2").Value = "SlugWatch"
$ws1.Range("H2").Value = "SlugWatch"
$ws1.Range("I2").Value = "SlugWatch"
$ws1.Range("J2").Value = "SlugWatch"
$ws1.Range("K2").Value = "SlugWatch"
$ws1.Range("M2").Value = "SlugWatch"

# --- SlugWatch2023!row2 (description) and row4 (purpose), per language column ---
# dk
$ws2.Range("C2").Value = 'SKADEVOLDER: Snegle er udbredt i mange havebrugs- og agermarker og kan skade afgrøder hele året. Med potentielt alvorlige konsekvenser for udbytte og/eller kvalitet er det vigtigt at styre bestande af nøgle-sneglearter. BESLUTNING: Behandling mod snegleangreb er mest effektiv, når den er rettet mod perioder med høj snegleaktivitet. SlugWatch-prognosen giver stedspecifikke snegleaktivitetsdata sammen med information om vejrforhold, der hjælper dig med at træffe informerede beslutninger om timingen og placeringen af ​​snegleovervågning og efterfølgende styring.
MODELLEN: SlugWatch-modellen forudsiger snegleaktivitet baseret på forudsagte vejrforhold.
PARAMETRENE: Modellen bruger jordtemperatur, lufttemperatur og nedbør.
KILDE: Skabt af CertisBelchim og udviklet af Farming Online Ltd. Valideret i Storbritannien og Frankrig og anses for egnet til forsigtig brug i Irland, Spanien, Belgien, Holland, Tyskland og Danmark.
ANTAGELSER: Denne model giver ingen forudsigelser om den faktiske forekomst af snegle, kun den relative risiko for snegleaktivitet, hvis den er til stede. Overvågning i marken er nødvendig for at vurdere den faktiske risiko for afgrøder i områder med høj snegleaktivitet. Overvej passende økonomiske tærskler før behandling.'
$ws2.Range("C4").Value = 'Se stedspecifikke live- og 7-dages prognosedata for snegleaktivitet'
# de
$ws2.Range("D2").Value = 'SCHADERREGER: Schnecken sind in vielen Gartenbau- und Ackerbaugebieten weit verbreitet und können das ganze Jahr über Ernten schädigen. Da die Auswirkungen auf den Ertrag und/oder die Qualität potenziell schwerwiegend sind, ist es wichtig, die Populationen der wichtigsten Schneckenarten zu kontrollieren. 
ENTSCHEIDUNG: Die Bekämpfung von Schneckenbefall ist am effektivsten, wenn sie in Zeiten hoher Schneckenaktivität erfolgt. Die SlugWatch-Vorhersage liefert ortsspezifische Daten zur Schneckenaktivität sowie Informationen zu den Wetterbedingungen und hilft so, fundierte Entscheidungen über Zeitpunkt und Ort der Schneckenüberwachung und des anschließenden Managements zu treffen. 
MODELL: Das SlugWatch-Modell sagt die Schneckenaktivität auf der Grundlage der vorhergesagten Wetterbedingungen voraus. 
PARAMETER: Das Modell verwendet Bodentemperatur, Lufttemperatur und Niederschlag. 
QUELLE: Erstellt von CertisBelchim und entwickelt von Farming Online Ltd. Es wurde im Vereinigten Königreich und in Frankreich validiert und gilt als geeignet für den Einsatz in Irland, Spanien, Belgien, den Niederlanden, Deutschland und Dänemark (unter Vorbehalt). 
ANNAHMEN: Dieses Modell macht keine Vorhersage über die tatsächliche Häufigkeit von Schnecken, sondern nur über das relative Risiko der Schneckenaktivität, wenn sie vorhanden sind. Um das tatsächliche Risiko einer hohen Schneckenaktivität für die Kulturpflanzen zu beurteilen, ist eine Überwachung im Feld erforderlich. Bitte beachten Sie die entsprechenden wirtschaftlichen Schwellenwerte für die Behandlung. '
$ws2.Range("D4").Value = 'Anzeige von ortsspezifischen Live- und 7-Tage-Vorhersagedaten zur Schneckenaktivität'
# nl
$ws2.Range("E2").Value = 'DE PLAAG: Slakken zijn op grote schaal aanwezig in veel tuinbouw- en akkerbouwgebieden en kunnen het hele jaar door schade toebrengen aan gewassen. Omdat de gevolgen voor opbrengst en/of kwaliteit ernstig kunnen zijn, is het belangrijk om de populaties van de belangrijkste slakkensoorten te beheersen. 
De BESLISSING: het beheer van slakkenplagen is het meest effectief wanneer het gericht is op periodes van hoge slakkenactiviteit. De SlugWatch prognose levert locatiespecifieke gegevens over de slakkenactiviteit, samen met informatie over de weersomstandigheden, en helpt bij het nemen van gefundeerde beslissingen over de timing en locatie van de slakkencontrole en het daaropvolgende beheer. 
HET MODEL: Het SlugWatch model voorspelt de slakkenactiviteit op basis van de voorspelde weersomstandigheden. 
DE PARAMETERS: Het model gebruikt bodemtemperatuur, luchttemperatuur en neerslag. 
BRON: Gemaakt door CertisBelchim en ontwikkeld door Farming Online Ltd. Gevalideerd in het VK en Frankrijk, geschikt geacht voor gebruik met de nodige terughoudendheid in Ierland, Spanje, België, Nederland, Duitsland en Denemarken. 
AANNAMES: Dit model doet geen voorspelling over de werkelijke slakkenpopulatie, alleen over het relatieve risico van slakkenactiviteit als ze aanwezig zijn. Om het werkelijke risico van een hoge slakkenactiviteit voor gewassen te beoordelen, is monitoring in het veld nodig. Raadpleeg de juiste economische drempelwaarden voor behandeling.
'
$ws2.Range("E4").Value = 'Bekijk live locatie-specifieke en 7-daagse prognoses over slakkenactiviteit.'
# sl
$ws2.Range("F2").Value = 'ŠKODLJIV ORGANIZEM: Polži, ki so pogosto prisotni na številnih vrtovih in poljih, lahko škodujejo pridelkom vse leto. Ker je vpliv na pridelek in/ali kakovost lahko resen, je pomembno upravljati populacije ključnih vrst polžev. ODLOČITEV: Upravljanje okužb s polži je najučinkovitejše, če je usmerjeno v obdobja velike aktivnosti polžev. Napoved SlugWatch zagotavlja podatke o aktivnosti polžev, specifične za posamezno lokacijo, skupaj z informacijami o vremenskih razmerah, kar pomaga pri sprejemanju informiranih odločitev o času in lokaciji spremljanja polžev ter poznejšem upravljanju. MODEL: Model SlugWatch napoveduje aktivnost polžev na podlagi napovedanih vremenskih razmer. PARAMETRI: Model uporablja temperaturo tal, temperaturo zraka in količino padavin. VIR: ustvarjeno v podjetju CertisBelchim, razvilo pa ga je podjetje Farming Online Ltd. Model je validiran v Združenem kraljestvu in Franciji, velja za primernega za previdno uporabo na Irskem, v Španiji, Belgiji, na Nizozemskem, v Nemčiji in na Danskem. PREDPOSTAVKE: Ta model ne predvideva dejanske številčnosti polžev, temveč le relativno tveganje aktivnosti polžev, če so prisotni. Za oceno dejanskega tveganja za posevke zaradi visoke aktivnosti polžev je potrebno spremljanje na terenu. Za obdelavo glejte ustrezne ekonomske pragove škodljivosti.'
$ws2.Range("F4").Value = 'Ogled podatkov o aktivnosti polžev v živo in sedemdnevne napovedi za posamezno lokacijo.'
# no
$ws2.Range("G2").Value = 'SKADEGJØREREN: Snegler er utbredt i mange hagebruk og dyrkbare marker, og kan skade avlinger hele året. Med potensielt alvorlige innvirkninger på avling og/eller kvalitet, er det viktig å håndtere populasjoner av sentrale sneglearter.
AVGJØRELSEN: Behandling mot snegleangrep er mest effektiv når den er rettet mot perioder med høy snegleaktivitet. SlugWatch-prognosen gir stedsspesifikke data for aktiviteten til snegler, sammen med informasjon om værforhold, og hjelper deg med å ta informerte beslutninger om timing og plassering av snegleovervåking og påfølgende håndtering.
MODELLEN: SlugWatch-modellen forutsier snegleaktivitet basert på varslede værforhold.
PARAMETRENE: Modellen bruker jordtemperatur, lufttemperatur og nedbør.
KILDE: Laget av CertisBelchim og utviklet av Farming Online Ltd. Validert i Storbritannia og Frankrike, og anses å være passende for forsiktig bruk i Irland, Spania, Belgia, Nederland, Tyskland og Danmark.
FORUTSETNINGER: Denne modellen gir ingen forutsigelser om den faktiske forekomsten av snegler, bare den relative risikoen for snegleaktivitet hvis de er tilstede. Overvåking i felt er nødvendig for å vurdere den faktiske risikoen for avlinger i områder med høy snegleaktivitet. Vurder hensiktsmessige økonomiske terskler før behandling.
'
$ws2.Range("G4").Value = 'Se lokasjons-spesifikke data for snegleaktivitet, live og som et 7-dagers varsel.'
# se
$ws2.Range("H2").Value = 'SKADEGÖRAREN: Sniglar finns i stor utsträckning i många trädgårds och åkermarker och kan skada grödor året runt. Sniglars effekter potentiellt allvarliga på avkastning och/eller kvalitet är det viktigt att hantera populationer av viktiga snigelarter.  BESLUTET: Hanteringen av snigelangrepp är mest effektiv när den riktas mot perioder med hög snigelaktivitet. SlugWatch-prognosen ger platsspecifika snigelaktivitetsdata, tillsammans med information om väderförhållanden, vilket hjälper till att fatta välgrundade beslut om tidpunkt och plats för snigelövervakning och efterföljande hantering.  MODELLEN: SlugWatch-modellen förutsäger snigelaktivitet baserat på prognostiserade väderförhållanden.  PARAMETRARNA: Modellen använder marktemperatur, lufttemperatur och nederbörd.  KÄLLA: Skapad av CertisBelchim och utvecklad av Farming Online Ltd. Validerad i Storbritannien och Frankrike, anses vara lämplig för användning med försiktighet i Irland, Spanien, Belgien, Nederländerna, Tyskland och Danmark.  ANTAGANDEN: Denna modell gör ingen förutsägelse om faktisk snigelmängd, endast den relativa risken för snigelaktivitet om de är närvarande. I fält krävs övervakning för att bedöma den faktiska risken för grödor med hög snigelaktivitet. Se lämpliga ekonomiska trösklar för behandling'
$ws2.Range("H4").Value = 'Se lokalspecifika data for snigelaktivitet, live och som en 7-dagarsprognos.'
# fi
$ws2.Range("I2").Value = 'KASVINTUHOOJA: Etanat esiintyvät laajalti monilla puutarhoilla ja pelloilla. Etanat voivat vahingoittaa viljelykasveja ympäri vuoden. On tärkeää hallita keskeisiä etanapopulaatioita niiden mahdollisesti vakavasti vaikuttaessa satoon ja/tai laatuun.
PÄÄTÖS: Etanoiden esiintymisen hallinta on tehokkainta, kun hallinta kohdennetaan etanoiden esiintymisen ajanjaksoihin. Slugwatch-ennuste tarjoaa paikkakohtaisesti etanoiden aktiivisuustiedot sekä sääolosuhteiden tiedot auttaen tekemään johtopäätöksiä etanan seurannan ja sitä seuraavan hallinnan ajoituksesta ja kohdentamisesta.
MALLI: Slugwatch -malli ennustaa etana-aktiivisuuden sääolosuhteiden ennustamiseen perustuen.
PARAMETRIT: Malli käyttää maaperän lämpötilaa, ilman lämpötilaa ja sadetta.
LÄHDE: Certisbelchimin luoma ja Iso -Britanniassa ja Ranskassa validoitu Farming Online Ltd., käytetään Irlannissa, Espanjassa, Belgiassa Alankomaissa, Saksassa ja Tanskassa.
OLETUKSET: Tämä malli ei tee ennustamista todellisesta etanoiden runsaudesta, ainoastaan suhteellisesta riskistä etanoiden aktiivisuudessa, silloin kun niitä esiintyy. Seurantaa tarvitaan todellisen riskin arvioimiseksi korkean etanan aktiivisuuden kasveilla. Pyydetään viittaamaan käsittelyn osalta taloudellisiin kynnysarvoihin.
'
$ws2.Range("I4").Value = 'Tarkastele sijaintikohtaisia reaaliaikaisia ja 7 päivän ennustetietoja etanan aktiivisuudesta'
# lt
$ws2.Range("J2").Value = 'KENKĖJAS: Šliužai, plačiai paplitę daugelyje sodininkystės ir daržininkystės pasėliuose ir gali kenkti pasėliams ištisus metus. Kadangi poveikis derliui ir (arba) kokybei potencialiai gali būti didelis, svarbu kontroliuoti pagrindinių rūšių šliužų populiacijas.
SPRENDIMAS: Šliužų plitimo kontrolė veiksmingiausia, kai ji vykdoma didelio šliužų aktyvumo laikotarpiais. „SlugWatch“ prognozė pateikia informaciją apie oro sąlygas ir duomenis apie šliužų aktyvumą konkrečioje vietovėje, padedant priimti pagrįstus sprendimus dėl šliužų stebėsenos ir tolimesnės kontrolės laiko ir vietos. 
MODELIS: „SlugWatch“ modelis prognozuoja šliužų aktyvumą pagal prognozuojamas oro sąlygas. 
PARAMETRAI: Modelis naudoja dirvožemio ir oro temperatūras ir kritulių kiekį. 
ŠALTINIS: Sukurtas „CertisBelchim“, išvystytas "Farming Online Ltd.". Patvirtintas Jungtinėje Karalystėje ir Prancūzijoje, laikomas tinkamu su atsargumu naudoti Airijoje, Ispanijoje, Belgijoje, Nyderlanduose, Vokietijoje ir Danijoje. 
PRIELAIDOS: Šis modelis nenumato tikrojo šliužų gausumo, o tik sąlyginę šliužų aktyvumo riziką, jei jie aptinkami. Norint įvertinti faktinį didelio šliužų aktyvumo pavojų pasėliams, būtina atlikti stebėseną pasėlyje. Remkitės atitinkamomis ekonominėmis žalingumo ribomis.'
$ws2.Range("J4").Value = 'Stebėkite konkrečios vietovės tiesioginės ir 7 dienų prognozių duomenis apie šliužų aktyvumą.'
# fr
$ws2.Range("K2").Value = 'LE RAVAGEUR : Largement présentes en grandes cultures et cultures légumières de plein champ, les limaces peuvent endommager les cultures tout au long de l''année. Les impacts sur le rendement et/ou la qualité pouvant être graves, il est important de gérer les populations des principales espèces de limaces. 
LA DÉCISION : La gestion des infestations de limaces est plus efficace lorsqu''elle est ciblée sur les périodes de forte activité des limaces. Les prévisions de SlugWatch fournissent des données sur l''activité des limaces en fonction du lieu, ainsi que des informations sur les conditions météorologiques, ce qui permet de prendre des décisions éclairées sur le moment et le lieu de la surveillance des limaces et de la gestion ultérieure. 
LE MODÈLE : Le modèle SlugWatch prévoit l''activité des limaces en fonction des conditions météorologiques prévues. 
LES PARAMÈTRES : Le modèle utilise la température du sol, la température de l''air et les précipitations. 
SOURCE : Créé par CertisBelchim et développé par Farming Online Ltd. Validé au Royaume-Uni et en France, considéré comme approprié pour une utilisation avec précaution en Irlande, Espagne, Belgique, Pays-Bas, Allemagne et Danemark. 
HYPOTHÈSES : Ce modèle ne fait aucune prédiction sur l''abondance réelle des limaces, seulement sur le risque relatif d''activité des limaces si elles sont présentes. Un suivi sur le terrain est nécessaire pour évaluer le risque réel pour les cultures d''une forte activité des limaces. Veuillez vous référer aux seuils économiques appropriés pour le traitement. '
$ws2.Range("K4").Value = 'Visualisez les données sur l''activité des limaces en direct et les prévisions à 7 jours, en fonction du lieu.'
# it
$ws2.Range("L2").Value = '"Il parassita: ampiamente presenti in molti campi di orticoltura e seminativi, le lumache possono danneggiare le colture durante tutto l''anno. Con impatti potenzialmente gravi sulla resa e/o sulla qualità, è importante gestire le popolazioni delle principali specie di lumache. 
LA DECISIONE: La gestione delle infestazioni di lumache è più efficace se mirata ai periodi di maggiore attività delle lumache. Le previsioni di SlugWatch forniscono dati specifici sull''attività delle limacce, insieme a informazioni sulle condizioni meteorologiche, aiutando a prendere decisioni informate sui tempi e i luoghi del monitoraggio delle limacce e della successiva gestione. 
IL MODELLO: il modello SlugWatch prevede l''attività delle lumache in base alle condizioni meteorologiche previste. 
I PARAMETRI: il modello utilizza la temperatura del suolo, la temperatura dell''aria e le precipitazioni. 
FONTE: Creato da CertisBelchim e sviluppato da Farming Online Ltd. Convalidato nel Regno Unito e in Francia, considerato appropriato per l''uso con cautela in Irlanda, Spagna, Belgio, Paesi Bassi, Germania e Danimarca. 
ASSUNZIONI: Questo modello non fa previsioni sull''effettiva abbondanza di lumache, ma solo sul rischio relativo di attività delle lumache in caso di loro presenza. È necessario un monitoraggio in campo per valutare il rischio effettivo per le colture di un''elevata attività delle limacce. Fare riferimento alle soglie economiche appropriate per il trattamento. "'
$ws2.Range("L4").Value = 'Visualizzare i dati sull''attività delle lumache in tempo reale e le previsioni per 7 giorni, specifiche per ogni località.'
# gr
$ws2.Range("M2").Value = ' Ο Εχθρός: Ευρέως παρόντες σε πολλά κηπευτικά και αροτραία χωράφια, οι γυμνοσάλιαγκες μπορούν να προκαλέσουν ζημιές στις καλλιέργειες όλο το χρόνο. Με δυνητικά σοβαρές επιπτώσεις στην απόδοση ή/και την ποιότητα, είναι σημαντικό να διαχειριστείτε τους πληθυσμούς των βασικών ειδών γυμνοσάλιαγκων. 
Η ΑΠΟΦΑΣΗ: Η διαχείριση των προσβολών από γυμνοσάλιαγκες είναι πιο αποτελεσματική όταν στοχεύει σε περιόδους υψηλής δραστηριότητας γυμνοσάλιαγκων. Η πρόβλεψη SlugWatch παρέχει δεδομένα για τη δραστηριότητα των γυμνοσάλιαγκων σε συγκεκριμένες τοποθεσίες, μαζί με πληροφορίες για τις καιρικές συνθήκες, βοηθώντας στη λήψη τεκμηριωμένων αποφάσεων σχετικά με το χρόνο και τη θέση παρακολούθησης των γυμνοσάλιαγκων και την επακόλουθη διαχείριση. 
ΤΟ ΜΟΝΤΕΛΟ: Το μοντέλο SlugWatch προβλέπει τη δραστηριότητα των γυμνοσάλιαγκων με βάση τις προβλεπόμενες καιρικές συνθήκες. 
ΟΙ ΠΑΡΑΜΕΤΡΟΙ: Το μοντέλο χρησιμοποιεί τη θερμοκρασία του εδάφους, τη θερμοκρασία του αέρα και τη βροχόπτωση. 
ΠΗΓΗ: Δημιουργήθηκε από την CertisBelchim και αναπτύχθηκε από την Farming Online Ltd. Επικυρώθηκε στο Ηνωμένο Βασίλειο και τη Γαλλία και θεωρείται κατάλληλο για χρήση με προσοχή στην Ιρλανδία, την Ισπανία, το Βέλγιο, τις Κάτω Χώρες, τη Γερμανία και τη Δανία. 
ΥΠΟΘΕΣΕΙΣ: Αυτό το μοντέλο δεν κάνει καμία πρόβλεψη σχετικά με την πραγματική αφθονία των γυμνοσάλιαγκων, παρά μόνο τον σχετικό κίνδυνο δραστηριότητας των γυμνοσάλιαγκων εάν υπάρχουν. Απαιτείται παρακολούθηση στον αγρό για να εκτιμηθεί ο πραγματικός κίνδυνος για τις καλλιέργειες από την υψηλή δραστηριότητα των γυμνοσάλιαγκων. Παρακαλούμε ανατρέξτε στα κατάλληλα οικονομικά κατώτατα όρια για την αντιμετώπιση. '
$ws2.Range("M4").Value = 'Προβολή δεδομένων δραστηριότητας γυμνοσάλιαγκων για συγκεκριμένη τοποθεσία σε ζωντανή μετάδοση και πρόβλεψη 7 ημερών'

# Row 2 holds the long multi-paragraph descriptions; writing them nudges
# Excel into giving the row an explicit wrapped-text height. Auto-fitting it
# drops that explicit height again so the row stays at its implicit default.
$ws2.Rows.Item(2).AutoFit()

# --- SlugWatch2023!row3 : model name ("SlugWatch") is unchanged across languages ---
$ws2.Range("C3").Value = "SlugWatch"
$ws2.Range("D3").Value = "SlugWatch"
$ws2.Range("F3").Value = "SlugWatch"
$ws2.Range("G3").Value = "SlugWatch"
$ws2.Range("H3").Value = "SlugWatch"
$ws2.Range("J3").Value = "SlugWatch"
$ws2.Range("K3").Value = "SlugWatch"

